$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 9093491
$ws.Range("I74").Value = 14288071
$ws.Range("J74").Value = 2976
$ws.Range("K74").Value = 14288071
$ws.Range("L74").Value = 2976
$ws.Range("M74").Value = -14287135
$ws.Range("N74").Value = -4848

# Row 77
$ws.Range("H77").Value = 9093491
$ws.Range("I77").Value = 14288071
$ws.Range("J77").Value = 2976
$ws.Range("K77").Value = 71440355
$ws.Range("L77").Value = 14880
$ws.Range("M77").Value = -71435675
$ws.Range("N77").Value = -24240

# Row 129
$ws.Range("H129").Value = 1090.5652
$ws.Range("I129").Value = 500
$ws.Range("J129").Value = 1117.409
$ws.Range("K129").Value = 1500
$ws.Range("L129").Value = 3352.227
$ws.Range("M129").Value = 3500
$ws.Range("N129").Value = -13352.227

# Row 137
$ws.Range("H137").Value = 4563.0244
$ws.Range("I137").Value = 4934.8696
$ws.Range("J137").Value = 4087.889
$ws.Range("K137").Value = 14804.6088
$ws.Range("L137").Value = 12263.667
$ws.Range("M137").Value = -12254.6088
$ws.Range("N137").Value = -17363.667

# Row 138
$ws.Range("H138").Value = 3457.36
$ws.Range("I138").Value = 1687.8889
$ws.Range("J138").Value = 3632.3625
$ws.Range("K138").Value = 5063.6667
$ws.Range("L138").Value = 10897.0875
$ws.Range("M138").Value = 76.33330000000024
$ws.Range("N138").Value = -21177.0875

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1320.125
$ws.Range("I2").Value = 1308.7142
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 1308.7142
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = -1195.7142

# Row 14
$ws.Range("H14").Value = 31949
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 31949
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 31949
$ws.Range("N14").Value = -32299

# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 32
$ws.Range("H32").Value = 14922.14
$ws.Range("I32").Value = 11322.5
$ws.Range("J32").Value = 18693.191
$ws.Range("K32").Value = 11322.5
$ws.Range("L32").Value = 18693.191
$ws.Range("M32").Value = -11035.5
$ws.Range("N32").Value = -19267.191

# Row 45
$ws.Range("H45").Value = 1181.7273
$ws.Range("I45").Value = 1199.9
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1199.9
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -822.9000000000001
$ws.Range("N45").Value = -1754

# Row 74
$ws.Range("H74").Value = 3611.7317
$ws.Range("I74").Value = 4187.115
$ws.Range("J74").Value = 2614.4
$ws.Range("K74").Value = 4187.115
$ws.Range("L74").Value = 2614.4
$ws.Range("M74").Value = -3313.115

# Row 77
$ws.Range("H77").Value = 3611.7317
$ws.Range("I77").Value = 4187.115
$ws.Range("J77").Value = 2614.4
$ws.Range("K77").Value = 20935.575
$ws.Range("L77").Value = 13072
$ws.Range("M77").Value = -16567.575

# Row 116
$ws.Range("H116").Value = 1320.125
$ws.Range("I116").Value = 1308.7142
$ws.Range("J116").Value = 1400
$ws.Range("K116").Value = 1308.7142
$ws.Range("L116").Value = 1400
$ws.Range("M116").Value = 985.2858000000001

# Row 122
$ws.Range("H122").Value = 4865
$ws.Range("I122").Value = 2063.3333
$ws.Range("J122").Value = 7666.6665
$ws.Range("K122").Value = 6189.999899999999
$ws.Range("L122").Value = 22999.9995
$ws.Range("M122").Value = -3739.999899999999

# Row 137
$ws.Range("H137").Value = 44186.332
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 44186.332
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 44186.332
$ws.Range("N137").Value = -54386.332

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1320.125
$ws.Range("I3").Value = 1308.7142
$ws.Range("J3").Value = 1400
$ws.Range("K3").Value = 1308.7142
$ws.Range("L3").Value = 1400
$ws.Range("M3").Value = -1194.7142

# Row 59
$ws.Range("H59").Value = 55390
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 55390
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 55390
$ws.Range("N59").Value = -57084

# Row 134
$ws.Range("H134").Value = 4400.1113
$ws.Range("I134").Value = 2236.6667
$ws.Range("J134").Value = 10890.444
$ws.Range("K134").Value = 6710.000100000001
$ws.Range("L134").Value = 32671.332
$ws.Range("M134").Value = -4175.000100000001
$ws.Range("N134").Value = -37741.33199999999

$ws = $wb.Worksheets.Item("CRP")
# Row 19
$ws.Range("H19").Value = 1126.6666
$ws.Range("I19").Value = 1126.6666
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1126.6666
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -956.6666

# Row 24
$ws.Range("H24").Value = 1126.6666
$ws.Range("I24").Value = 1126.6666
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1126.6666
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -956.6666

# Row 31
$ws.Range("H31").Value = 5043.757
$ws.Range("I31").Value = 1704.3
$ws.Range("J31").Value = 6280.593
$ws.Range("K31").Value = 1704.3
$ws.Range("L31").Value = 6280.593
$ws.Range("M31").Value = -1409.3
$ws.Range("N31").Value = -6870.593

# Row 34
$ws.Range("H34").Value = 5043.757
$ws.Range("I34").Value = 1704.3
$ws.Range("J34").Value = 6280.593
$ws.Range("K34").Value = 1704.3
$ws.Range("L34").Value = 6280.593
$ws.Range("M34").Value = -1502.3
$ws.Range("N34").Value = -6684.593

# Row 38
$ws.Range("H38").Value = 24000
$ws.Range("I38").Value = 10000
$ws.Range("J38").Value = 27500
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 27500
$ws.Range("M38").Value = -9623
$ws.Range("N38").Value = -28254

# Row 46
$ws.Range("H46").Value = 24000
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 27500
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 27500
$ws.Range("M46").Value = -9789
$ws.Range("N46").Value = -27922

# Row 94
$ws.Range("H94").Value = 1675.8572
$ws.Range("I94").Value = 1594.9
$ws.Range("J94").Value = 1749.4546
$ws.Range("K94").Value = 1594.9
$ws.Range("L94").Value = 1749.4546
$ws.Range("M94").Value = -1143.9
$ws.Range("N94").Value = -2651.4546

$ws = $wb.Worksheets.Item("CUL")
# Row 19
$ws.Range("H19").Value = 3000
$ws.Range("I19").Value = 3000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 9000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -8826

# Row 23
$ws.Range("H23").Value = 241.42857
$ws.Range("I23").Value = 119.666664
$ws.Range("J23").Value = 274.63635
$ws.Range("K23").Value = 358.999992
$ws.Range("L23").Value = 823.90905
$ws.Range("M23").Value = -123.999992
$ws.Range("N23").Value = -1293.90905

# Row 113
$ws.Range("H113").Value = 639.7222
$ws.Range("I113").Value = 640.28
$ws.Range("J113").Value = 638.4545000000001
$ws.Range("K113").Value = 1920.84
$ws.Range("L113").Value = 1915.3635
$ws.Range("M113").Value = 249.1600000000001
$ws.Range("N113").Value = -6255.3635

# Row 117
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("M117").ClearContents()
$ws.Range("N117").ClearContents()

# Row 121
$ws.Range("H121").Value = 1859.4789
$ws.Range("I121").Value = 489.8
$ws.Range("J121").Value = 1963.2424
$ws.Range("K121").Value = 1469.4
$ws.Range("L121").Value = 5889.7272
$ws.Range("M121").Value = -159.4000000000001
$ws.Range("N121").Value = -8509.727200000001

# Row 129
$ws.Range("H129").Value = 3572.9167
$ws.Range("I129").Value = 2752.3076
$ws.Range("J129").Value = 4542.727
$ws.Range("K129").Value = 8256.9228
$ws.Range("L129").Value = 13628.181
$ws.Range("M129").Value = -3256.9228
$ws.Range("N129").Value = -23628.181

# Row 131
$ws.Range("H131").Value = 11910693
$ws.Range("I131").Value = 55579276
$ws.Range("J131").Value = 1079.3636
$ws.Range("K131").Value = 166737828
$ws.Range("L131").Value = 3238.0908
$ws.Range("M131").Value = -166732788
$ws.Range("N131").Value = -13318.0908

$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()

# Row 8
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()

# Row 137
$ws.Range("H137").Value = 57665.11
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 57665.11
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 57665.11
$ws.Range("N137").Value = -67865.11

$ws = $wb.Worksheets.Item("LTW")
# Row 24
$ws.Range("H24").Value = 19999
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 19999
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 19999
$ws.Range("N24").Value = -20685

# Row 132
$ws.Range("H132").Value = 5065.766
$ws.Range("I132").Value = 2359.5715
$ws.Range("J132").Value = 9053.842000000001
$ws.Range("K132").Value = 7078.7145
$ws.Range("L132").Value = 27161.526
$ws.Range("M132").Value = -4548.7145
$ws.Range("N132").Value = -32221.526

$ws = $wb.Worksheets.Item("WVR")
# Row 23
$ws.Range("H23").Value = 60508
$ws.Range("I23").Value = 21005
$ws.Range("J23").Value = 100011
$ws.Range("K23").Value = 21005
$ws.Range("L23").Value = 100011
$ws.Range("M23").Value = -20776
$ws.Range("N23").Value = -100469

# Row 124
$ws.Range("H124").Value = 24426
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 24426
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 24426
$ws.Range("N124").Value = -34246
